$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(106, 1).Value = 0.000075866666666666704
$ws.Cells.Item(106, 2).Value = 0.000017733333333333301
$ws.Cells.Item(106, 3).Value = 0.000029733333333333302
$ws.Cells.Item(107, 1).Value = -0.00018573333333333299
$ws.Cells.Item(107, 2).Value = 0.000017600000000000001
$ws.Cells.Item(107, 3).Value = 0.000052666666666666702
$ws.Cells.Item(108, 1).Value = 0.000019733333333333299
$ws.Cells.Item(108, 2).Value = 0.0000134666666666667
$ws.Cells.Item(108, 3).Value = 0.00012
$ws.Cells.Item(109, 1).Value = 0.000024133333333333301
$ws.Cells.Item(109, 2).Value = 0.000034533333333333303
$ws.Cells.Item(109, 3).Value = 0.000046400000000000003
$ws.Cells.Item(110, 1).Value = 0.000023200000000000001
$ws.Cells.Item(110, 2).Value = 0.000049466666666666699
$ws.Cells.Item(110, 3).Value = 0.0001076
$ws.Cells.Item(111, 1).Value = 0.000035066666666666701
$ws.Cells.Item(111, 2).Value = 0.0000341333333333333
$ws.Cells.Item(111, 3).Value = 0.00017026666666666699
$ws.Cells.Item(112, 1).Value = -0.0000238666666666667
$ws.Cells.Item(112, 2).Value = 0.000014399999999999999
$ws.Cells.Item(112, 3).Value = 0.000080799999999999999
$ws.Cells.Item(113, 1).Value = 0.000024533333333333301
$ws.Cells.Item(113, 2).Value = 0.000036000000000000001
$ws.Cells.Item(113, 3).Value = 0.000204133333333333
$ws.Cells.Item(114, 1).Value = 0.000099866666666666704
$ws.Cells.Item(114, 2).Value = 0.000053733333333333302
$ws.Cells.Item(114, 3).Value = 0.000158533333333333
$ws.Cells.Item(115, 1).Value = -0.000000000000000000000052222404641385099
$ws.Cells.Item(115, 2).Value = 0.000024000000000000001
$ws.Cells.Item(115, 3).Value = 0.000070533333333333304
$ws.Cells.Item(116, 1).Value = 0.00000026666666666666699
$ws.Cells.Item(116, 2).Value = 0.000023200000000000001
$ws.Cells.Item(116, 3).Value = 0.000088800000000000004
$ws.Cells.Item(117, 1).Value = 0.0000071999999999999997
$ws.Cells.Item(117, 2).Value = 0.000037866666666666702
$ws.Cells.Item(117, 3).Value = 0.00018866666666666701
$ws.Cells.Item(118, 1).Value = 0.000018133333333333301
$ws.Cells.Item(118, 2).Value = 0.0000396
$ws.Cells.Item(118, 3).Value = 0.00019133333333333299
$ws.Cells.Item(119, 1).Value = -0.000062799999999999995
$ws.Cells.Item(119, 2).Value = 0.0000030666666666666702
$ws.Cells.Item(119, 3).Value = 0.00016906666666666699
$ws.Cells.Item(120, 1).Value = -0.000111866666666667
$ws.Cells.Item(120, 2).Value = 0.0000156
$ws.Cells.Item(120, 3).Value = 0.000210666666666667
$ws.Cells.Item(121, 1).Value = -0.00000253333333333333
$ws.Cells.Item(121, 2).Value = 0.000032933333333333298
$ws.Cells.Item(121, 3).Value = 0.00028266666666666701
$ws.Cells.Item(122, 1).Value = -0.000016533333333333299
$ws.Cells.Item(122, 2).Value = 0.000033866666666666699
$ws.Cells.Item(122, 3).Value = 0.00017866666666666701
$ws.Cells.Item(123, 1).Value = -0.0000076000000000000001
$ws.Cells.Item(123, 2).Value = 0.000011733333333333299
$ws.Cells.Item(123, 3).Value = 0.000116533333333333
$ws.Cells.Item(124, 1).Value = 0.0000079999999999999996
$ws.Cells.Item(124, 2).Value = 0.000031466666666666702
$ws.Cells.Item(124, 3).Value = 0.00026453333333333301
$ws.Cells.Item(125, 1).Value = -0.0000433333333333333
$ws.Cells.Item(125, 2).Value = -0.0000034666666666666702
$ws.Cells.Item(125, 3).Value = 0.00012653333333333301
$ws.Cells.Item(126, 1).Value = 0.000014933333333333299
$ws.Cells.Item(126, 2).Value = 0.000019599999999999999
$ws.Cells.Item(126, 3).Value = 0.00022253333333333299
$ws.Cells.Item(127, 1).Value = -0.000033333333333333301
$ws.Cells.Item(127, 2).Value = 0.0000184
$ws.Cells.Item(127, 3).Value = 0.00022893333333333301
$ws.Cells.Item(128, 1).Value = -0.000042933333333333297
$ws.Cells.Item(128, 2).Value = 0.0000221333333333333
$ws.Cells.Item(128, 3).Value = 0.000090533333333333302
$ws.Cells.Item(129, 1).Value = -0.000029066666666666701
$ws.Cells.Item(129, 2).Value = 0.000022266666666666699
$ws.Cells.Item(129, 3).Value = 0.000177466666666667
$ws.Cells.Item(130, 1).Value = -0.000023733333333333301
$ws.Cells.Item(130, 2).Value = 0.000012799999999999999
$ws.Cells.Item(130, 3).Value = 0.000072266666666666698
$ws.Cells.Item(131, 1).Value = -0.000056933333333333299
$ws.Cells.Item(131, 2).Value = -0.000013333333333333299
$ws.Cells.Item(131, 3).Value = 0.000045599999999999997
$ws.Cells.Item(132, 1).Value = -0.00021493333333333299
$ws.Cells.Item(132, 2).Value = -0.000031466666666666702
$ws.Cells.Item(132, 3).Value = 0.000020800000000000001
$ws.Cells.Item(133, 1).Value = 0.0000081333333333333306
$ws.Cells.Item(133, 2).Value = 0.0000092
$ws.Cells.Item(133, 3).Value = 0.00015880000000000001
$ws.Cells.Item(134, 1).Value = -0.000032266666666666701
$ws.Cells.Item(134, 2).Value = 0.0000047999999999999998
$ws.Cells.Item(134, 3).Value = 0.00013426666666666701
$ws.Cells.Item(135, 1).Value = -0.0000244
$ws.Cells.Item(135, 2).Value = 0.00000013333333333333299
$ws.Cells.Item(135, 3).Value = 0.00018546666666666701
$ws.Cells.Item(136, 1).Value = -0.000097999999999999997
$ws.Cells.Item(136, 2).Value = -0.000020000000000000002
$ws.Cells.Item(136, 3).Value = 0.00014573333333333299
$ws.Cells.Item(137, 1).Value = -0.000019066666666666699
$ws.Cells.Item(137, 2).Value = 0.000015866666666666699
$ws.Cells.Item(137, 3).Value = 0.00012973333333333301
$ws.Cells.Item(138, 1).Value = 0.000032933333333333298
$ws.Cells.Item(138, 2).Value = 0.000012133333333333301
$ws.Cells.Item(138, 3).Value = 0.00019760000000000001
$ws.Cells.Item(139, 1).Value = 0.0000106666666666667
$ws.Cells.Item(139, 2).Value = -0.0000019999999999999999
$ws.Cells.Item(139, 3).Value = 0.00022693333333333301
$ws.Cells.Item(140, 1).Value = -0.0000145333333333333
$ws.Cells.Item(140, 2).Value = -0.000012
$ws.Cells.Item(140, 3).Value = 0.00021346666666666701
$ws.Cells.Item(141, 1).Value = -0.000067199999999999994
$ws.Cells.Item(141, 2).Value = -0.000025466666666666698
$ws.Cells.Item(141, 3).Value = 0.00010066666666666699
$ws.Cells.Item(142, 1).Value = -0.000020533333333333302
$ws.Cells.Item(142, 2).Value = 0.0000094666666666666703
$ws.Cells.Item(142, 3).Value = 0.00018560000000000001
$ws.Cells.Item(143, 1).Value = -0.0000088000000000000004
$ws.Cells.Item(143, 2).Value = 0.0000011999999999999999
$ws.Cells.Item(143, 3).Value = 0.00023599999999999999
$ws.Cells.Item(144, 1).Value = 0.0000341333333333333
$ws.Cells.Item(144, 2).Value = -0.0000058666666666666701
$ws.Cells.Item(144, 3).Value = 0.000065733333333333296
$ws.Cells.Item(145, 1).Value = -0.000047200000000000002
$ws.Cells.Item(145, 2).Value = -0.000039733333333333301
$ws.Cells.Item(145, 3).Value = 0.000047066666666666702
$ws.Cells.Item(146, 1).Value = -0.000038266666666666698
$ws.Cells.Item(146, 2).Value = -0.000019066666666666699
$ws.Cells.Item(146, 3).Value = 0.0000266666666666667
$ws.Cells.Item(147, 1).Value = -0.000041333333333333299
$ws.Cells.Item(147, 2).Value = -0.000021333333333333301
$ws.Cells.Item(147, 3).Value = 0.000185333333333333
$ws.Cells.Item(148, 1).Value = -0.0000524
$ws.Cells.Item(148, 2).Value = -0.000015866666666666699
$ws.Cells.Item(148, 3).Value = 0.000075599999999999994
$ws.Cells.Item(149, 1).Value = -0.000118266666666667
$ws.Cells.Item(149, 2).Value = -0.0000097333333333333305
$ws.Cells.Item(149, 3).Value = 0.000071333333333333296
$ws.Cells.Item(150, 1).Value = -0.0000486666666666667
$ws.Cells.Item(150, 2).Value = -0.000027333333333333301
$ws.Cells.Item(150, 3).Value = 0.00012400000000000001
$ws.Cells.Item(151, 1).Value = -0.000023733333333333301
$ws.Cells.Item(151, 2).Value = 0.000020533333333333302
$ws.Cells.Item(151, 3).Value = 0.00011
$ws.Cells.Item(152, 1).Value = -0.0000104
$ws.Cells.Item(152, 2).Value = 0.0000086666666666666695
$ws.Cells.Item(152, 3).Value = 0.00030733333333333302
$ws.Cells.Item(153, 1).Value = 0.000033733333333333297
$ws.Cells.Item(153, 2).Value = 0.000046799999999999999
$ws.Cells.Item(153, 3).Value = 0.00027253333333333298
